$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G20:G29").HorizontalAlignment = -4131
$ws.Range("G32").HorizontalAlignment = -4131
